$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 319; this shifts existing rows 319..467 down to 320..468
$ws.Rows.Item(319).EntireRow.Insert()

# Populate the newly inserted row 319 with the new weekly record.
# Values mirror the former row 319 record, except for the date (D),
# minimum price (K), weighted average price (M) and Price $/Kg (P).
$ws.Range("A319").Value = 9
$ws.Range("B319").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C319").Value = "Metropolitana"
$ws.Range("D319").Value = 45134
$ws.Range("E319").Value = 13
$ws.Range("F319").Value = 300000001
$ws.Range("G319").Value = "Rabanito"
$ws.Range("H319").Value = "Sin especificar"
$ws.Range("I319").Value = "Primera"
$ws.Range("J319").Value = 7000
$ws.Range("K319").Value = 3000
$ws.Range("L319").Value = 4000
$ws.Range("M319").Value = 3500
$ws.Range("N319").Value = '$/cien unidades (volumen en unidades)'
$ws.Range("O319").Value = "Provincia de Chacabuco"
$ws.Range("P319").Value = 35
$ws.Range("Q319").Value = 100
$ws.Range("R319").Value = "Hortaliza"

Write-Host "Done"
